$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 is a brand-new row (DemoCSR / crmsfa / Paypal / Priya / Paul).
# Write the A3 cell first, then the C/D/E cells of row 2 (whose text
# changes from CTS/Rohith/Nandakumar to SalesForce/Jane/Hendrix), then
# fill in the rest of row 3, and finally B3 (which reuses the existing
# "crmsfa" shared string) -- this ordering reproduces the shared-string
# table layout seen in the target workbook.
$ws.Range("A3").Value2 = "DemoCSR"

$ws.Range("C2").Value2 = "SalesForce"
$ws.Range("C3").Value2 = "Paypal"

$ws.Range("D2").Value2 = "Jane"
$ws.Range("E2").Value2 = "Hendrix"

$ws.Range("D3").Value2 = "Priya"
$ws.Range("E3").Value2 = "Paul"

$ws.Range("B3").Value2 = "crmsfa"

# Columns D:E got narrower once the new (shorter) names replaced the old
# ones -- mirror that with an explicit width tweak.
$ws.Columns.Item(4).ColumnWidth = 9.666666666666666
$ws.Columns.Item(5).ColumnWidth = 9.333333333333334

# Selection moved to F3 after the edits.
$ws.Range("F3").Select()
